$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Word keeps exactly one "_GoBack" bookmark, re-stamped at the
#    location of the most recent edit. Drop it from its old home
#    (around the picture in the "Conditional Codes" paragraph) before
#    we make the new edit below, so it ends up only in the new spot.
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ------------------------------------------------------------------
# 2. Insert a new bold run "Data Movement Instructions. " right
#    before the "Page 205" run (its own run, not merged into the
#    neighbour, carrying the w:hint="eastAsia" attribute), followed
#    immediately by the relocated "_GoBack" bookmark.
# ------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("Page 205", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null

$insPos = $target.Start
$insertRange = $d.Range($insPos, $insPos)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:b/></w:rPr><w:t xml:space="preserve">Data Movement Instructions. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$insertRange.InsertXML($xml)
